$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.323.94"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.863.42"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'243.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'0.7006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("D8").Value = "'0.07906"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "'24.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").Value = "'0.07812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.15%  "
$ws.Range("D12").Value = "1.872.57"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "'5.142"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "'92.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "'0.6977"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "'6.534"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "'0.000008539"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "29.349.83"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'248.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "2.123.63"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").Value = "'12.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "'7.586"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").Value = "'8.963"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").Value = "'160.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "'18.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").Value = "'1.583"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.16%  "
$ws.Range("D30").Value = "'4.299"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").Value = "'4.241"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'1.207"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'0.05236"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "'1.886"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("D35").Value = "'0.7569"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'1.177"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'2.698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "1.275.62"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'0.01865"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "'0.9026"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "'110.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D43").Value = "'5.956"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.95%  "
$ws.Range("D44").Value = "'70.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.48%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "2.022.16"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("D48").Value = "'9.618"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "'1.789"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'0.4281"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.27%  "
